# Changing Dependency refset to be of type CidCidStr.
# Insert a new "name (String)" column (B) into the FHIMDependenciesRS sheet,
# shifting the existing client/supplier/Notes columns one position to the
# right, and give every existing data row an empty-string placeholder value
# in the new column (matches the `""` placeholder convention already used
# elsewhere in this workbook, e.g. FHIMAttributesRS "name (String)" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FHIMDependenciesRS")

# Insert a new blank column at B; this shifts old B/C/D -> C/D/E and moves
# the column width metadata along with them.
$ws.Columns("B").Insert()

# Approximate the intended column width for the new "name (String)" column
# (closest value reachable given the host's column-width rounding).
$ws.Columns("B").ColumnWidth = 17

# New header cell.
$ws.Range("B1").Value = "name (String)"

# New data cells: literal placeholder value `""` for every populated row.
$ws.Range("B2").Value = """"""
$ws.Range("B3").Value = """"""
$ws.Range("B4").Value = """"""
$ws.Range("B5").Value = """"""
$ws.Range("B7").Value = """"""

# Match the author's final selection (cell B7, the new "name" cell of the
# last data row).
$ws.Activate()
$ws.Range("B7").Select()
